$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new column G with header "Tipo" and values for each certificate row
$ws.Range("G1").Value = "Tipo"
$ws.Range("G2").Value = "Primario"
$ws.Range("G3").Value = "Secundario"
$ws.Range("G4").Value = "Vencido"

# Update the active selection to match the recorded end-state (H7)
$ws.Range("H7").Select()
